# Adds a new acronym-key entry for "SIOM" (Standard Input Output Matrix) to the
# "Key to Variables" sheet, in the "io-model" Top Level Folder group, just above
# the existing "URPbIC" row, with "high" importance to update for a new country.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new row right above the current row 161 ("URPbIC"); this shifts
# row 161 and everything below it down by one row and, by default, copies
# the formatting of the row above (row 160).
$ws.Rows.Item(161).Insert()

# Row 160's F-column fill ("medium") was copied into the new row by the
# insert above; re-stamp F161's format from an existing "high"-rated row
# (now at row 179, after the shift) so the fill color matches other "high"
# importance cells.
$ws.Range("F179").Copy()
$ws.Range("F161").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row's content.
$ws.Range("A161").Value = "io-model"
$ws.Range("B161").Value = "SIOM"
$ws.Range("C161").Value = "Standard Input Output Matrix"
$ws.Range("F161").Value = "high"
